$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45188 -> 45189) for every data row (rows 2 through 330).
$ws.Range("C2:C330").Value = 45189
